$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H6").Value = 667298.75
$ws.Range("I6").Value = 1428838.8
$ws.Range("J6").Value = 951.25
$ws.Range("K6").Value = 4286516.4
$ws.Range("L6").Value = 2853.75
$ws.Range("M6").Value = -4286404.4
$ws.Range("N6").Value = -3077.75
$ws.Range("H18").Value = 1161.4615
$ws.Range("I18").Value = 758.25
$ws.Range("K18").Value = 758.25
$ws.Range("M18").Value = -474.25
$ws.Range("H43").Value = 2431.2222
$ws.Range("I43").Value = 2996.8333
$ws.Range("J43").Value = 1300
$ws.Range("K43").Value = 2996.8333
$ws.Range("L43").Value = 1300
$ws.Range("M43").Value = -2927.8333
$ws.Range("N43").Value = -1438
$ws.Range("H74").Value = 3039.72
$ws.Range("I74").Value = 3000.15
$ws.Range("J74").Value = 3198
$ws.Range("K74").Value = 3000.15
$ws.Range("L74").Value = 3198
$ws.Range("M74").Value = -2064.15
$ws.Range("N74").Value = -5070
$ws.Range("H77").Value = 3039.72
$ws.Range("I77").Value = 3000.15
$ws.Range("J77").Value = 3198
$ws.Range("K77").Value = 15000.75
$ws.Range("L77").Value = 15990
$ws.Range("M77").Value = -10320.75
$ws.Range("N77").Value = -25350

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H8").Value = 0
$ws.Range("J8").Value = 0
$ws.Range("L8").Value = 0
$ws.Range("N8").Value = $null
$ws.Range("H32").Value = 10715.907
$ws.Range("I32").Value = 8731.530000000001
$ws.Range("J32").Value = 30162.8
$ws.Range("K32").Value = 8731.530000000001
$ws.Range("L32").Value = 30162.8
$ws.Range("M32").Value = -8444.530000000001
$ws.Range("N32").Value = -30736.8
$ws.Range("H61").Value = 1910.6333
$ws.Range("I61").Value = 1606
$ws.Range("J61").Value = 2367.5833
$ws.Range("K61").Value = 1606
$ws.Range("L61").Value = 2367.5833
$ws.Range("M61").Value = -1394
$ws.Range("N61").Value = -2791.5833
$ws.Range("H74").Value = 1196.3871
$ws.Range("I74").Value = 1196.4615
$ws.Range("J74").Value = 1196
$ws.Range("K74").Value = 1196.4615
$ws.Range("L74").Value = 1196
$ws.Range("M74").Value = -322.4614999999999
$ws.Range("N74").Value = -2944
$ws.Range("H77").Value = 1196.3871
$ws.Range("I77").Value = 1196.4615
$ws.Range("J77").Value = 1196
$ws.Range("K77").Value = 5982.307499999999
$ws.Range("L77").Value = 5980
$ws.Range("M77").Value = -1614.307499999999
$ws.Range("N77").Value = -14716
$ws.Range("H88").Value = 1725
$ws.Range("I88").Value = 1500
$ws.Range("J88").Value = 1800
$ws.Range("K88").Value = 1500
$ws.Range("L88").Value = 1800
$ws.Range("M88").Value = -1094
$ws.Range("N88").Value = -2612
$ws.Range("H91").Value = 1725
$ws.Range("I91").Value = 1500
$ws.Range("J91").Value = 1800
$ws.Range("K91").Value = 1500
$ws.Range("L91").Value = 1800
$ws.Range("M91").Value = -96
$ws.Range("N91").Value = -4608
$ws.Range("H132").Value = 3966.6
$ws.Range("I132").Value = 3326.5
$ws.Range("J132").Value = 6527
$ws.Range("K132").Value = 9979.5
$ws.Range("L132").Value = 19581
$ws.Range("M132").Value = -7449.5
$ws.Range("N132").Value = -24641
$ws.Range("H136").Value = 1910.6333
$ws.Range("I136").Value = 1606
$ws.Range("J136").Value = 2367.5833
$ws.Range("K136").Value = 4818
$ws.Range("L136").Value = 7102.749899999999
$ws.Range("M136").Value = -2268
$ws.Range("N136").Value = -12202.7499

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H12").Value = 1995
$ws.Range("I12").Value = 1995
$ws.Range("K12").Value = 1995
$ws.Range("M12").Value = -1827
$ws.Range("H42").Value = 398000
$ws.Range("J42").Value = 398000
$ws.Range("L42").Value = 398000
$ws.Range("N42").Value = -398656
$ws.Range("H86").Value = 86463.234
$ws.Range("I86").Value = 101658.37
$ws.Range("J86").Value = 2890
$ws.Range("K86").Value = 101658.37
$ws.Range("L86").Value = 2890
$ws.Range("M86").Value = -100535.37
$ws.Range("N86").Value = -5136
$ws.Range("H89").Value = 86463.234
$ws.Range("I89").Value = 101658.37
$ws.Range("J89").Value = 2890
$ws.Range("K89").Value = 508291.85
$ws.Range("L89").Value = 14450
$ws.Range("M89").Value = -502675.85
$ws.Range("N89").Value = -25682
$ws.Range("H107").Value = 333488770
$ws.Range("J107").Value = 0
$ws.Range("L107").Value = 0
$ws.Range("N107").Value = $null
$ws.Range("H134").Value = 2470.2727
$ws.Range("I134").Value = 2170.2144
$ws.Range("J134").Value = 4150.6
$ws.Range("K134").Value = 6510.6432
$ws.Range("L134").Value = 12451.8
$ws.Range("M134").Value = -3975.6432
$ws.Range("N134").Value = -17521.8

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H2").Value = 4
$ws.Range("I2").Value = 4
$ws.Range("K2").Value = 4
$ws.Range("M2").Value = 109
$ws.Range("H31").Value = 2660.4187
$ws.Range("I31").Value = 852.2857
$ws.Range("J31").Value = 3533.3103
$ws.Range("K31").Value = 852.2857
$ws.Range("L31").Value = 3533.3103
$ws.Range("M31").Value = -557.2857
$ws.Range("N31").Value = -4123.3103
$ws.Range("H34").Value = 2660.4187
$ws.Range("I34").Value = 852.2857
$ws.Range("J34").Value = 3533.3103
$ws.Range("K34").Value = 852.2857
$ws.Range("L34").Value = 3533.3103
$ws.Range("M34").Value = -650.2857
$ws.Range("N34").Value = -3937.3103
$ws.Range("H58").Value = 2387.5264
$ws.Range("I58").Value = 2466.5833
$ws.Range("J58").Value = 2252
$ws.Range("K58").Value = 2466.5833
$ws.Range("L58").Value = 2252
$ws.Range("M58").Value = -2263.5833
$ws.Range("N58").Value = -2658
$ws.Range("H59").Value = 24666.666
$ws.Range("J59").Value = 25000
$ws.Range("L59").Value = 25000
$ws.Range("N59").Value = -27290
$ws.Range("H107").Value = 684.6667
$ws.Range("I107").Value = 691.875
$ws.Range("J107").Value = 670.25
$ws.Range("K107").Value = 691.875
$ws.Range("L107").Value = 670.25
$ws.Range("M107").Value = 1228.125
$ws.Range("N107").Value = -4510.25
$ws.Range("H136").Value = 2387.5264
$ws.Range("I136").Value = 2466.5833
$ws.Range("J136").Value = 2252
$ws.Range("K136").Value = 7399.749899999999
$ws.Range("L136").Value = 6756
$ws.Range("M136").Value = -4849.749899999999
$ws.Range("N136").Value = -11856

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 1154.5897
$ws.Range("I5").Value = 1346
$ws.Range("J5").Value = 1079.3928
$ws.Range("K5").Value = 4038
$ws.Range("L5").Value = 3238.1784
$ws.Range("M5").Value = -3926
$ws.Range("N5").Value = -3462.1784
$ws.Range("H12").Value = 43.36842
$ws.Range("I12").Value = 19.2
$ws.Range("J12").Value = 52
$ws.Range("K12").Value = 57.59999999999999
$ws.Range("L12").Value = 156
$ws.Range("M12").Value = 115.4
$ws.Range("N12").Value = -502
$ws.Range("H34").Value = 1054.5454
$ws.Range("J34").Value = 1222.2222
$ws.Range("L34").Value = 3666.6666
$ws.Range("N34").Value = -3834.6666
$ws.Range("H135").Value = 1154.5897
$ws.Range("I135").Value = 1346
$ws.Range("J135").Value = 1079.3928
$ws.Range("K135").Value = 12114
$ws.Range("L135").Value = 9714.5352
$ws.Range("M135").Value = -9579
$ws.Range("N135").Value = -14784.5352

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H5").Value = 4500
$ws.Range("I5").Value = 0
$ws.Range("K5").Value = 0
$ws.Range("M5").Value = $null
$ws.Range("H52").Value = 7816.25
$ws.Range("J52").Value = 9166.666999999999
$ws.Range("L52").Value = 9166.666999999999
$ws.Range("N52").Value = -9684.666999999999
$ws.Range("H80").Value = 91004500
$ws.Range("I80").Value = 250258750
$ws.Range("J80").Value = 2064.2856
$ws.Range("K80").Value = 250258750
$ws.Range("L80").Value = 2064.2856
$ws.Range("M80").Value = -250257752
$ws.Range("N80").Value = -4060.2856
$ws.Range("H83").Value = 91004500
$ws.Range("I83").Value = 250258750
$ws.Range("J83").Value = 2064.2856
$ws.Range("K83").Value = 1251293750
$ws.Range("L83").Value = 10321.428
$ws.Range("M83").Value = -1251288758
$ws.Range("N83").Value = -20305.428

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H136").Value = 1392.3556
$ws.Range("I136").Value = 1135.2059
$ws.Range("J136").Value = 2187.182
$ws.Range("K136").Value = 3405.6177
$ws.Range("L136").Value = 6561.545999999999
$ws.Range("M136").Value = -855.6176999999998
$ws.Range("N136").Value = -11661.546

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H45").Value = 10450
$ws.Range("J45").Value = 10450
$ws.Range("L45").Value = 10450
$ws.Range("N45").Value = -11432
